# Update automatico via Actualizar 06-01-2020 02-38-27
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's data as row 80
$newRow = 80
$ws.Cells.Item($newRow, 1).Value = 43982
$ws.Cells.Item($newRow, 2).Value = 404
$ws.Cells.Item($newRow, 3).Value = 108
$ws.Cells.Item($newRow, 4).Value = 488
$ws.Cells.Item($newRow, 5).Value = 22
$ws.Cells.Item($newRow, 6).Value = 40

# Match formatting of the row above (date format on column A, plain/centered numeric elsewhere)
$ws.Range("A79").Copy()
$ws.Range("A80").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B79:F79").Copy()
$ws.Range("B80:F80").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Grow the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F80"))

# Move selection to the newly added cell, matching the saved view state
$ws.Range("F80").Select()
